$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Planner (high prio)" was split across 3 runs (with spell-check
#    proofErr markers around "prio"). Collapse it back into a single run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Planner (high prio)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Planner (high prio)", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Same pattern for the Profile bullet (split around "btn").
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Profile – Remove follow button from your profile view… replace with “Upload nudes” btn.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Profile – Remove follow button from your profile view… replace with “Upload nudes” btn.", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) "Find wanderers and " / "Pending wanderers activities" (with the
#    _GoBack bookmark sitting between the two runs) becomes one merged
#    run, and a brand new list item "Chat activity color changes" is
#    added right after it (inheriting the same ListParagraph / numId=3 /
#    bold / green formatting), with the _GoBack bookmark now living at
#    the end of that new paragraph's text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Find wanderers and Pending wanderers activities", $true, $false, $false, $false, $false,
    $true, 1, $false, "Find wanderers and Pending wanderers activities", 2) | Out-Null

$targetIdx = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Find wanderers and Pending wanderers activities*") {
        $targetIdx = $idx
    }
}

$target = $d.Paragraphs.Item($targetIdx)
$cr = [char]13

# Append a paragraph break plus the new sentence, with a throw-away
# trailing marker character. Inserting a new paragraph this way (rather
# than InsertParagraphAfter) makes the split paragraph naturally inherit
# the source paragraph's pPr/rPr (style, numbering, bold, color).
$target.Range.InsertAfter($cr + "Chat activity color changesX")

$newPara = $d.Paragraphs.Item($targetIdx + 1)

# Position right before the marker character 'X' (i.e. right after the
# real text, still inside the paragraph, not sitting on the paragraph
# mark itself - bookmarking exactly on the pilcrow silently no-ops).
$bmPos = $newPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Remove the throw-away marker character now that the bookmark is anchored.
$xRange = $d.Range($bmPos, $bmPos + 1)
$xRange.Text = ""
